# Actualización desde MV -datos-
# Append two new daily rows (06-09-2021 and 07-09-2021) to the single sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Date = "06-09-2021"; B = 1.42; C = 2.4;  D = 2.75; E = 2.4; F = -1.02 },
    @{ Date = "07-09-2021"; B = 1.44; C = 2.27; D = 2.89; E = 2.4; F = -0.98 }
)

$firstNewRow = 171
$lastNewRow = $firstNewRow + $newRows.Count - 1

# Column A holds a dd-mm-yyyy date label that must stay plain text (it
# mirrors every earlier "Serie" entry, which are shared strings, not real
# dates). Forcing the cells to Text before assigning keeps Excel's automatic
# date recognition from turning the label into a date serial, and
# ClearFormats() afterwards drops the temporary "@" number format again so
# the cells are left on the default (unstyled) format, same as every other
# data row in the sheet.
$dateRange = $ws.Range("A$firstNewRow`:A$lastNewRow")
$dateRange.NumberFormat = "@"

$lastRow = $firstNewRow - 1
foreach ($entry in $newRows) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $entry.Date
    $ws.Cells.Item($lastRow, 2).Value = $entry.B
    $ws.Cells.Item($lastRow, 3).Value = $entry.C
    $ws.Cells.Item($lastRow, 4).Value = $entry.D
    $ws.Cells.Item($lastRow, 5).Value = $entry.E
    $ws.Cells.Item($lastRow, 6).Value = $entry.F
}

$dateRange.ClearFormats()
